# Update error-message cells on the "Training Results" sheet so that the
# embedded onnx node-name counters reflect the new run's numbering.
# (Diff only touches text content of column C cells; no structural changes.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

$ws.Range("C16").Value = "C:\Users\COCO\onnxruntime_training_cuda_python\orttraining\orttraining\python\orttraining_pybind_state.cc:621 onnxruntime::python::addObjectMethodsForTraining::<lambda_6dd399ad6691adab5d0e0423ed8ce22d>::operator () [ONNXRuntimeError] : 1 : FAIL : Type Error: Type parameter (T) of Optype (Sub) bound to different types (tensor(float) and tensor(double) in node (onnx::Pow::27433_Grad/Sub_1).`n"

$ws.Range("C39").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27552): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C51").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27642): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C52").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27644): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C64").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27697): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C65").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27699): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C69").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27717): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C70").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27719): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C93").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::27815): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C102").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::27852): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C223").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::28302): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C239").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::28307): X typestr: T, has unsupported type: tensor(uint8)"
